# Rows 3, 4 and 5 of the "Artfynd" sheet get their observation data cyclically
# rotated: the data that used to live in row 5 moves up to row 3, the data
# that used to live in row 3 moves down to row 4, and the data that used to
# live in row 4 moves down to row 5. The columns below are the only ones
# whose values actually differ between the three rows; every other column
# (C, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AX, AY) is already
# identical across rows 3-5, so there is nothing to move for those.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","D","E","F","G","H","I","J","K","M","N","Q","R","AF","AJ","AK","AO")

# Column I ("Antal") holds its number as text (e.g. "1") rather than a true
# numeric value in this sheet, so it needs to be written back with a leading
# apostrophe to keep it stored as text instead of being auto-converted to a
# number. Columns I, J, K, N and AF can also legitimately hold an explicit
# empty-text value (as opposed to a truly blank/empty cell); writing a plain
# "" normally collapses to a blank cell, so those also need the leading
# apostrophe (a lone "'" keeps the cell text-typed with empty content).
$textCols = @("I","J","K","N","AF")

# Snapshot the current (pre-edit) contents of rows 3, 4 and 5 for the columns
# that need to move, before anything gets overwritten.
$row3 = @{}
$row4 = @{}
$row5 = @{}
foreach ($col in $cols) {
    $row3[$col] = $ws.Range($col + "3").Value2
    $row4[$col] = $ws.Range($col + "4").Value2
    $row5[$col] = $ws.Range($col + "5").Value2
}

function Convert-ForWrite($col, $val) {
    if (@($textCols) -contains $col) {
        if ($val -eq $null) {
            return $null
        } else {
            return "'" + $val
        }
    }
    return $val
}

# Apply the rotation: row5 -> row3, row3 -> row4, row4 -> row5.
foreach ($col in $cols) {
    $v3 = Convert-ForWrite $col $row5[$col]
    $v4 = Convert-ForWrite $col $row3[$col]
    $v5 = Convert-ForWrite $col $row4[$col]

    $ws.Range($col + "3").Value = $v3
    $ws.Range($col + "4").Value = $v4
    $ws.Range($col + "5").Value = $v5
}
